$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 40 - 2020-04-23 (serial 43944)
$ws.Cells.Item(40, 1).Value = 43944
$ws.Cells.Item(40, 2).Value = -763
$ws.Cells.Item(40, 4).Value = -118
$ws.Cells.Item(40, 6).Value = 422
$ws.Cells.Item(40, 7).Value = 1296

# Row 41 - 2020-04-24 (serial 43945)
$ws.Cells.Item(41, 1).Value = 43945
$ws.Cells.Item(41, 2).Value = -754
$ws.Cells.Item(41, 4).Value = -108
$ws.Cells.Item(41, 6).Value = 437
$ws.Cells.Item(41, 7).Value = 1184

# Match date style of column A used in existing rows by copying the
# existing formatting (e.g. from A39) rather than re-declaring a format,
# so we don't create a duplicate/new number-format style entry.
$ws.Range("A39").Copy()
$ws.Range("A40:A41").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the selected/active cell to reflect the new last-entry position
$ws.Range("F42").Select()
